# Update countries dataset (new case numbers + resulting sort-by-"Casos totales" order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered list: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
  @("Estados Unidos",99116,13681,2465,95125,2463,231,1526),
  @("Italia",86498,5909,10950,66414,3732,919,9134),
  @("China",81340,55,74588,3460,1034,5,3292),
  @("España",64059,6273,9357,49768,4165,569,4934),
  @("Alemania",50178,6240,6658,43182,23,71,338),
  @("Francia",32964,3809,5700,25269,3787,299,1995),
  @("Iran",32332,2926,11133,18821,2893,144,2378),
  @("Reino Unido",14543,2885,135,13649,163,181,759),
  @("Suiza",12928,1117,1530,11167,203,39,231),
  @("Corea del Sur",9332,91,4528,4665,59,8,139),
  @("Paises Bajos",8603,1172,3,8054,761,112,546),
  @("Austria",7657,748,225,7374,128,9,58),
  @("Belgica",7284,1049,858,6137,690,69,289),
  @("Turquia",5698,2069,42,5564,241,17,92),
  @("Canada",4610,567,228,4329,120,14,53),
  @("Portugal",4268,724,43,4149,71,16,76),
  @("Noruega",3717,345,6,3692,70,5,19),
  @("Brasil",3417,432,6,3319,296,15,92),
  @("Australia",3180,130,170,2997,23,0,13),
  @("Suecia",3046,206,16,2938,214,15,92),
  @("Israel",3035,342,79,2944,49,4,12),
  @("Chequia",2279,354,11,2259,34,0,9),
  @("Malasia",2161,130,259,1876,54,3,26),
  @("Irlanda",2121,302,5,2094,47,3,22),
  @("Dinamarca",2046,169,1,1993,109,11,52),
  @("Chile",1610,304,43,1562,7,1,5),
  @("Luxemburgo",1605,152,40,1550,25,6,15),
  @("Ecuador",1595,192,3,1556,58,2,36),
  @("Polonia",1389,168,7,1366,3,0,16),
  @("Japon",1387,0,359,981,57,0,47),
  @("Pakistan",1331,130,23,1298,7,1,10),
  @("Rumania",1292,263,115,1151,32,3,26),
  @("Sudafrica",1170,243,12,1157,7,1,1),
  @("Tailandia",1136,91,97,1034,11,1,5),
  @("Arabia Saudita",1104,92,35,1066,6,0,3),
  @("Indonesia",1046,153,46,913,0,9,87),
  @("Finlandia",1041,83,10,1024,32,2,7),
  @("Rusia",1036,196,45,987,8,1,4),
  @("Grecia",966,74,52,886,66,1,28),
  @("Islandia",890,88,97,791,18,0,2),
  @("India",887,160,73,794,0,0,20),
  @("Filipinas",803,96,31,718,1,9,54),
  @("Singapur",732,49,183,547,17,0,2),
  @("Crucero",712,0,597,105,15,0,10),
  @("Panama",674,0,2,663,20,0,9),
  @("Peru",635,55,14,612,14,0,9),
  @("Eslovenia",632,70,10,613,14,3,9),
  @("Argentina",589,0,72,504,0,1,13),
  @("Croacia",586,91,37,546,14,0,3),
  @("Mexico",585,110,4,573,1,2,8),
  @("Republica Dominicana",581,93,3,558,0,10,20),
  @("Estonia",575,37,11,563,7,0,1),
  @("Catar",562,13,43,519,6,0,0),
  @("Colombia",539,48,8,525,0,0,6),
  @("Egipto",536,41,116,390,0,6,30),
  @("Serbia",528,71,42,478,25,1,8),
  @("Hong Kong",518,64,111,403,5,0,4),
  @("Barein",466,8,227,235,1,0,4),
  @("Irak",458,76,122,296,0,4,40),
  @("Argelia",409,42,29,354,0,1,26),
  @("Emiratos Arabes Unidos",405,72,52,351,2,0,2),
  @("Libano",391,23,27,356,3,2,8),
  @("Nueva Zelanda",368,0,37,331,1,0,0),
  @("Lituania",358,59,1,352,2,1,5),
  @("Marruecos",333,58,11,301,1,10,21),
  @("Armenia",329,39,28,300,6,0,1),
  @("Hungria",300,39,34,256,6,0,10),
  @("Bulgaria",293,29,9,281,8,0,3),
  @("Letonia",280,36,1,279,0,0,0),
  @("Eslovaquia",269,43,2,267,1,0,0),
  @("Principado de Andorra",267,43,1,263,11,0,3),
  @("Taiwan",267,15,30,235,0,0,2),
  @("Costa Rica",263,32,2,259,5,0,2),
  @("Uruguay",238,0,0,238,3,0,0),
  @("Jordania",235,23,17,218,0,0,0),
  @("Bosnia y Herzegovina",232,41,5,223,1,1,4),
  @("Tunez",227,30,2,219,10,1,6),
  @("Ucrania",226,30,5,216,0,0,5),
  @("Kuwait",225,17,57,168,11,0,0),
  @("San Marino",223,15,4,198,12,0,21),
  @("Republica de Macedonia",219,18,3,213,1,0,3),
  @("Moldavia",199,22,2,195,33,1,2),
  @("Albania",186,12,31,147,3,2,8),
  @("Burkina Faso",180,28,12,159,0,2,9),
  @("Azerbaiyan",165,43,15,147,6,0,3),
  @("Vietnam",163,10,20,143,3,0,0),
  @("Republica de Chipre",162,16,4,153,3,2,5),
  @("Kazajistan",150,37,3,146,0,0,1),
  @("Reunion",145,10,1,144,0,0,0),
  @("Islas Feroe",144,4,54,90,2,0,0),
  @("Malta",139,5,2,137,1,0,0),
  @("Ghana",136,4,1,131,1,0,4),
  @("Oman",131,22,23,108,0,0,0),
  @("Senegal",119,14,11,108,0,0,0),
  @("Brunei",115,1,11,104,1,0,0),
  @("Afganistan",110,16,2,104,0,0,4),
  @("Venezuela",107,0,31,75,2,0,1),
  @("Sri Lanka",106,0,7,99,5,0,0),
  @("Costa de Marfil",101,5,3,98,0,0,0),
  @("Camboya",99,1,11,88,1,0,0),
  @("Mauricio",94,13,0,92,1,0,2),
  @("Bielorrusia",94,8,32,62,2,0,0),
  @("Estado de Palestina",91,5,17,73,0,0,1),
  @("Camerun",88,13,2,84,0,1,2),
  @("Uzbekistan",88,13,5,82,8,1,1),
  @("Georgia",83,4,14,69,1,0,0),
  @("Martinica",81,0,0,80,12,0,1),
  @("Cuba",80,13,4,74,2,0,2),
  @("Montenegro",75,6,0,74,1,0,1),
  @("Guadalupe",73,0,0,72,4,0,1),
  @("Nigeria",70,5,3,66,0,0,1),
  @("Honduras",68,1,0,67,0,0,1),
  @("Trinidad yTobago",66,1,1,63,0,1,2),
  @("Bolivia",61,0,0,61,0,0,0),
  @("Kirguistan",58,14,0,58,0,0,0),
  @("Liechtenstein",56,0,0,56,0,0,0),
  @("Gibraltar",55,20,14,41,0,0,0),
  @("Ruanda",54,4,0,54,0,0,0),
  @("Paraguay",52,11,1,48,1,0,3),
  @("Consejo Danes para los Refugiados",51,0,2,46,0,0,3),
  @("Mayotte",50,14,0,50,0,0,0),
  @("Banglades",48,4,11,32,1,0,5),
  @("Monaco",42,9,1,41,0,0,0),
  @("Puerto Rico",39,0,1,36,0,0,2),
  @("Macao",34,1,10,24,0,0,0),
  @("Guam",32,0,0,31,0,0,1),
  @("Kenia",31,0,1,29,2,0,1),
  @("Polinesia Francesa",30,0,0,30,0,0,0),
  @("Isla de Man",29,3,0,29,0,0,0),
  @("Aruba",28,0,1,27,0,0,0),
  @("Guatemala",28,3,4,23,1,0,1),
  @("Guayana Francesa",28,0,6,22,0,0,0),
  @("Madagascar",26,3,0,26,0,0,0),
  @("Jamaica",26,0,2,23,0,0,1),
  @("Togo",25,1,1,24,0,0,0),
  @("Barbados",24,0,0,24,0,0,0),
  @("Zambia",22,6,0,22,0,0,0),
  @("Uganda",18,4,0,18,0,0,0),
  @("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
  @("Etiopia",16,4,0,16,0,0,0),
  @("Maldivas",16,3,9,7,0,0,0),
  @("Nueva Caledonia",15,1,0,15,0,0,0),
  @("Bermudas",15,0,2,13,0,0,0),
  @("El Salvador",13,0,0,13,0,0,0),
  @("Tanzania",13,0,1,12,0,0,0),
  @("Guinea Ecuatorial",12,0,0,12,0,0,0),
  @("Republica de Yibuti",12,1,0,12,0,0,0),
  @("Dominica",11,0,0,11,0,0,0),
  @("Mali",11,7,0,11,0,0,0),
  @("San Martin (Parte Francesa)",11,0,0,11,0,0,0),
  @("Mongolia",11,0,0,11,0,0,0),
  @("Niger",10,0,0,9,0,0,1),
  @("Groenlandia",10,4,2,8,0,0,0),
  @("Suazilandia",9,3,0,9,0,0,0),
  @("Bahamas",9,0,1,8,0,0,0),
  @("Birmania",8,3,0,8,0,0,0),
  @("Haiti",8,0,0,8,0,0,0),
  @("Guinea",8,4,0,8,0,0,0),
  @("Surinam",8,0,0,8,0,0,0),
  @("Islas Caimanes",8,0,0,7,0,0,1),
  @("Namibia",8,0,2,6,0,0,0),
  @("Granada",7,0,0,7,0,0,0),
  @("Seychelles",7,0,0,7,0,0,0),
  @("Mozambique",7,0,0,7,0,0,0),
  @("Antigua y Barbuda",7,0,0,7,0,0,0),
  @("Gabon",7,0,0,6,0,0,1),
  @("Curazao",7,0,2,4,0,0,1),
  @("Laos",6,0,0,6,0,0,0),
  @("Benin",6,0,0,6,0,0,0),
  @("Eritrea",6,0,0,6,0,0,0),
  @("Montserrat",5,0,0,5,0,0,0),
  @("San Bartolome",5,2,0,5,0,0,0),
  @("Siria",5,0,0,5,0,0,0),
  @("Fiyi",5,0,0,5,0,0,0),
  @("Cabo Verde",5,0,0,4,0,0,1),
  @("Guyana",5,0,0,4,0,0,1),
  @("Zimbabue",5,2,0,4,0,0,1),
  @("Angola",4,0,0,4,0,0,0),
  @("Congo",4,0,0,4,0,0,0),
  @("Santa Sede",4,0,0,4,0,0,0),
  @("Nepal",4,1,1,3,0,0,0),
  @("Republica del Chad",3,0,0,3,0,0,0),
  @("Butan",3,1,0,3,0,0,0),
  @("Mauritania",3,0,0,3,0,0,0),
  @("Republica de Africa Central",3,0,0,3,0,0,0),
  @("San Martin (Parte Holandesa)",3,0,0,3,0,0,0),
  @("Somalia",3,1,0,3,0,0,0),
  @("Liberia",3,0,0,3,0,0,0),
  @("Sudan",3,0,0,2,0,0,1),
  @("Santa Lucia",3,0,1,2,0,0,0),
  @("Gambia",3,0,0,2,0,0,1),
  @("Anguila",2,0,0,2,0,0,0),
  @("Islas Turcas y Caicos",2,0,0,2,0,0,0),
  @("Belice",2,0,0,2,0,0,0),
  @("Guinea-Bisau",2,0,0,2,0,0,0),
  @("San Cristobal y Nieves",2,0,0,2,0,0,0),
  @("Islas Virgenes Britanicas",2,0,0,2,0,0,0),
  @("Nicaragua",2,0,0,1,0,1,1),
  @("San Vicente y las Granadinas",1,0,0,1,0,0,0),
  @("Libia",1,0,0,1,0,0,0),
  @("Timor Oriental",1,0,0,1,0,0,0),
  @("Papua Nueva Guinea",1,0,0,1,0,0,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rec = $data[$i]
  $ws.Cells.Item($r, 1).Value = $rec[0]
  $ws.Cells.Item($r, 2).Value = $rec[1]
  $ws.Cells.Item($r, 3).Value = $rec[2]
  $ws.Cells.Item($r, 4).Value = $rec[3]
  $ws.Cells.Item($r, 5).Value = $rec[4]
  $ws.Cells.Item($r, 6).Value = $rec[5]
  $ws.Cells.Item($r, 7).Value = $rec[6]
  $ws.Cells.Item($r, 8).Value = $rec[7]
}

Write-Host "Updated $($data.Count) country rows"
